$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.464105367660522
$ws.Range("B1").Value = 1.697013854980469
$ws.Range("C1").Value = 2.156684637069702
$ws.Range("D1").Value = 3.470783472061157
$ws.Range("E1").Value = 4.040699005126953
